# Refresh the 2021.08.07 "Change"/running-total columns that were re-pulled
# from the source (new totals are a few votes higher than the stale numbers).
# Each "Change" cell is set to its new value, and the adjacent running-total
# cell is updated to Total = PreviousTotal + Change, matching the sheet's own
# convention (seen e.g. in the live formula D2 = E2-C2 on the MMV sheet).
$wb = $excel.ActiveWorkbook

# --- AMV: Anime Movies (Votes) -- row 32 refresh for 2021.08.07 ---
$ws = $wb.Worksheets.Item("AMV")
$ws.Range("R32").Value = 134401
$ws.Range("S32").Value = 1590816

# --- AFV: Anime Favorites (Votes) -- 2021.08.07 Change/Total columns (D/E) ---
$ws = $wb.Worksheets.Item("AFV")
$ws.Range("D2").Value = -12067
$ws.Range("E2").Value = 182184
$ws.Range("D3").Value = -5111
$ws.Range("E3").Value = 156845
$ws.Range("D4").Value = -1753
$ws.Range("E4").Value = 154697
$ws.Range("D5").Value = -19028
$ws.Range("E5").Value = 136964
$ws.Range("D6").Value = -14210
$ws.Range("E6").Value = 131656
$ws.Range("D7").Value = -10712
$ws.Range("E7").Value = 128920
$ws.Range("D9").Value = -12175
$ws.Range("E9").Value = 80951
$ws.Range("D10").Value = -16337
$ws.Range("E10").Value = 61913
$ws.Range("D11").Value = -902
$ws.Range("E11").Value = 76053
$ws.Range("D12").Value = -3147
$ws.Range("E12").Value = 72658
$ws.Range("D13").Value = -2788
$ws.Range("E13").Value = 71712
$ws.Range("D14").Value = -12975
$ws.Range("E14").Value = 59071
$ws.Range("D15").Value = -3180
$ws.Range("E15").Value = 67114
$ws.Range("D17").Value = -4226
$ws.Range("E17").Value = 64089
$ws.Range("D18").Value = -5879
$ws.Range("E18").Value = 62055
$ws.Range("D19").Value = -2988
$ws.Range("E19").Value = 64775
$ws.Range("D21").Value = -17272
$ws.Range("E21").Value = 46153
$ws.Range("D22").Value = -2165
$ws.Range("E22").Value = 57918
$ws.Range("D23").Value = -7653
$ws.Range("E23").Value = 49731
$ws.Range("D24").Value = -5343
$ws.Range("E24").Value = 51896
$ws.Range("D25").Value = -2186
$ws.Range("E25").Value = 55004
$ws.Range("D26").Value = -5597
$ws.Range("E26").Value = 48786
$ws.Range("D28").Value = -1308
$ws.Range("E28").Value = 49188
$ws.Range("D30").Value = -6879
$ws.Range("E30").Value = 42097
$ws.Range("D31").Value = -6856
$ws.Range("E31").Value = 41180
$ws.Range("D32").Value = -5806
$ws.Range("E32").Value = 42109
$ws.Range("D34").Value = -3477
$ws.Range("E34").Value = 43716
$ws.Range("D35").Value = 4103
$ws.Range("E35").Value = 49720
$ws.Range("D36").Value = -520
$ws.Range("E36").Value = 42669
$ws.Range("D37").Value = -1266
$ws.Range("E37").Value = 38826
$ws.Range("D38").Value = 1420
$ws.Range("E38").Value = 40829
$ws.Range("D41").Value = -3530
$ws.Range("E41").Value = 32635
$ws.Range("D42").Value = -4954
$ws.Range("E42").Value = 29511
$ws.Range("D43").Value = 669
$ws.Range("E43").Value = 34750
$ws.Range("D44").Value = -4728
$ws.Range("E44").Value = 29080
$ws.Range("D49").Value = -1543
$ws.Range("E49").Value = 29914
$ws.Range("D51").Value = 322
$ws.Range("E51").Value = 31186

# --- MRV: Manga Ranking (Votes) -- row 41 refresh for 2021.08.07 ---
$ws = $wb.Worksheets.Item("MRV")
$ws.Range("P41").Value = -0.01000000000000156
$ws.Range("Q41").Value = 8.699999999999999

# --- MMV: Manga Movies... (Votes) -- 2021.08.07 Change/Total columns (L/M) ---
$ws = $wb.Worksheets.Item("MMV")
$ws.Range("L2").Value = 48125
$ws.Range("M2").Value = 498664
$ws.Range("L3").Value = 29901
$ws.Range("M3").Value = 410353
$ws.Range("L4").Value = 12588
$ws.Range("M4").Value = 350544
$ws.Range("L5").Value = 42809
$ws.Range("M5").Value = 427677
$ws.Range("L6").Value = 28542
$ws.Range("M6").Value = 370469
$ws.Range("L7").Value = 25311
$ws.Range("M7").Value = 329477
$ws.Range("L8").Value = 18364
$ws.Range("M8").Value = 315796
$ws.Range("L9").Value = 11842
$ws.Range("M9").Value = 288863
$ws.Range("L10").Value = 30700
$ws.Range("M10").Value = 332500
$ws.Range("L11").Value = 5723
$ws.Range("M11").Value = 242430
$ws.Range("L12").Value = 11096
$ws.Range("M12").Value = 236889
$ws.Range("L13").Value = 27816
$ws.Range("M13").Value = 276162
$ws.Range("L14").Value = 10792
$ws.Range("M14").Value = 219281
$ws.Range("L15").Value = 19573
$ws.Range("M15").Value = 259863
$ws.Range("L16").Value = 14742
$ws.Range("M16").Value = 210455
$ws.Range("L17").Value = 13172
$ws.Range("M17").Value = 195951
$ws.Range("L18").Value = 21248
$ws.Range("M18").Value = 234466
$ws.Range("L19").Value = 6354
$ws.Range("M19").Value = 172136
$ws.Range("L20").Value = 4384
$ws.Range("M20").Value = 160426
$ws.Range("L21").Value = 8776
$ws.Range("M21").Value = 170991
$ws.Range("L22").Value = 7789
$ws.Range("M22").Value = 167646
$ws.Range("L23").Value = 23941
$ws.Range("M23").Value = 211219
$ws.Range("L24").Value = 43411
$ws.Range("M24").Value = 245332
$ws.Range("L25").Value = 6956
$ws.Range("M25").Value = 155178
$ws.Range("L26").Value = 6487
$ws.Range("M26").Value = 147934
$ws.Range("L27").Value = 6018
$ws.Range("M27").Value = 145219
$ws.Range("L28").Value = 6358
$ws.Range("M28").Value = 148406
$ws.Range("L29").Value = 19625
$ws.Range("M29").Value = 181123
$ws.Range("L30").Value = 8231
$ws.Range("M30").Value = 145534
$ws.Range("L31").Value = 5844
$ws.Range("M31").Value = 136222
$ws.Range("L32").Value = 7816
$ws.Range("M32").Value = 140884
$ws.Range("L33").Value = 4981
$ws.Range("M33").Value = 132477
$ws.Range("L34").Value = 4949
$ws.Range("M34").Value = 132158
$ws.Range("L36").Value = 6725
$ws.Range("M36").Value = 134077
$ws.Range("L37").Value = 9843
$ws.Range("M37").Value = 143650
$ws.Range("L38").Value = 19835
$ws.Range("M38").Value = 174610
$ws.Range("L39").Value = 50870
$ws.Range("M39").Value = 254723
$ws.Range("L40").Value = 12061
$ws.Range("M40").Value = 146405
$ws.Range("L41").Value = 14104
$ws.Range("M41").Value = 148686
$ws.Range("L42").Value = 7153
$ws.Range("M42").Value = 129371
$ws.Range("L43").Value = 16582
$ws.Range("M43").Value = 160644
$ws.Range("L45").Value = 16113
$ws.Range("M45").Value = 147017
$ws.Range("L47").Value = 22533
$ws.Range("M47").Value = 161600
$ws.Range("L51").Value = 15238
$ws.Range("M51").Value = 149635
$ws.Range("L52").Value = 90933
$ws.Range("M52").Value = 258705
$ws.Range("L53").Value = 13790
$ws.Range("M53").Value = 145226
$ws.Range("L54").Value = 58819
$ws.Range("M54").Value = 185690
$ws.Range("L55").Value = 15788
$ws.Range("M55").Value = 137382
$ws.Range("M56").Value = 137797
